$d = $word.ActiveDocument
$cr = [char]13

# The "Contact Phone Number" paragraph lives in the single cell of the
# first table. Paragraphs.Item()/.Last are unreliable on this host, so
# walk the .Next() chain from .First looking for the paragraph's text
# instead of relying on a fixed index.
$cell = $d.Tables.Item(1).Cell(1, 1)
$contactPara = $cell.Range.Paragraphs.First
while (($contactPara -ne $null) -and ($contactPara.Range.Text -notlike "Contact Phone Number*")) {
    $contactPara = $contactPara.Next()
}

# Insert a brand-new paragraph right before "Contact Phone Number: ..."
# that reads "IRB #: 201901813" (inherits the same paragraph formatting
# because InsertBefore splits off of contactPara's own paragraph mark).
$contactPara.Range.InsertBefore("IRB #: 201901813" + $cr)

# Re-seat the document's "_GoBack" bookmark (Word tracks the most recent
# edit point with it) to sit between "IRB #" and ": 201901813" - exactly
# where the cursor would be after typing "IRB #" and before typing the
# rest. Adding a bookmark with an existing name moves it, removing the
# old one automatically.
$goBackRange = $d.Content.Duplicate
$goBackRange.Find.Execute("IRB #") | Out-Null
$goBackRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goBackRange)
